# Update the Correspond Handoff/Handback DateTime values for the
# 5eec1f23-... / fed47d42-... handback rows on the zh-cn and de-de sheets.
# This reflects a re-run of the handback report generation, which produced
# newer handoff/handback timestamps for that row pair.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 23 & 24, column D = Correspond Handoff Datetime,
# column G = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D23").Value = "2016-03-01 09:51:07"
$wsZhCn.Range("D24").Value = "2016-03-01 09:51:07"
$wsZhCn.Range("G23").Value = "2016-03-01 09:52:04"
$wsZhCn.Range("G24").Value = "2016-03-01 09:52:04"

# de-de sheet: same rows/columns
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D23").Value = "2016-03-01 09:51:18"
$wsDeDe.Range("D24").Value = "2016-03-01 09:51:18"
$wsDeDe.Range("G23").Value = "2016-03-01 09:52:22"
$wsDeDe.Range("G24").Value = "2016-03-01 09:52:22"
